$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.405.77"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "3.992.01"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.67"
$ws.Range("E5").Value = "  +6.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.37"
$ws.Range("E6").Value = "  +14.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.689"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.801"
$ws.Range("E9").Value = "  +4.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.189"
$ws.Range("E10").Value = "  +9.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.67"
$ws.Range("E11").Value = "  +6.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000340"
$ws.Range("E12").Value = "  +3.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.73"
$ws.Range("E13").Value = "  +4.33%  "
$ws.Range("D14").Value = "4.632.06"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "3.988.42"
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.41"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "21.01"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "73.320.40"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "464.93"
$ws.Range("E21").Value = "  +4.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.85"
$ws.Range("E22").Value = "  +6.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.74"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("E24").Value = "  -4.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.36"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.24"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.27"
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.72"
$ws.Range("E28").Value = "  -3.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.96"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.58"
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.03"
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "14.06"
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "49.71"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0000104"
$ws.Range("E34").Value = "  +14.99%  "
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.44"
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "641.97"
$ws.Range("E37").Value = "  -6.56%  "
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.149"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.42"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0488"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.24"
$ws.Range("E44").Value = "  +37.18%  "
$ws.Range("E45").Value = "  -5.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000308"
$ws.Range("E46").Value = "  +11.73%  "
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.42"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.63"
$ws.Range("E49").Value = "  -4.79%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.01"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.822.20"
$ws.Range("E51").Value = "  +1.43%  "
